# Apply the table style change described by the commit's diff:
#   Slide 5 -> the table shape's <a:tableStyleId> changes from
#   {C6B27CAC-9E35-4557-A2C5-46A9FE308BFE} to
#   {1A9630F7-A04C-4EDE-992B-0D4AB31BFD3D}
#
# Table.Style is a read-only reflection of the applied style; PowerPoint's
# object model requires Table.ApplyStyle(styleId) to actually change it.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTable) {
        $sh.Table.ApplyStyle("{1A9630F7-A04C-4EDE-992B-0D4AB31BFD3D}")
    }
}
